# Update the PCB BOM workbook:
#  1. Rename the sheet tab date from 2025-10-14 to 2025-10-15
#  2. Update the LED component row (row 4) with the new part info

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet (tab name embeds the date)
$ws.Name = "BOM_Board1_1_PCB_2025-10-15"

# 2. Update the LED row (row 4: No. = 3) with new component data
$ws.Range("C4").Value = "YLED1206R"
$ws.Range("D4").Value = "LED"
$ws.Range("E4").Value = "LED1206-FD"
$ws.Range("G4").Value = "YLED1206R"
$ws.Range("H4").Value = "YONGYUTAI(永裕泰)"
$ws.Range("I4").Value = "C28310439"
